$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1098
$ws.Range("H3").Value = 8348
$ws.Range("J3").Value = 1178
$ws.Range("J4").Value = 265
$ws.Range("J5").Value = 88
$ws.Range("J6").Value = 1582
$ws.Range("H7").Value = 26000
$ws.Range("J7").Value = 4211

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 40
$ws.Range("J6").Value = 42
$ws.Range("I7").Value = 818
$ws.Range("J7").Value = 123
$ws.Range("J8").Value = 268
$ws.Range("J11").Value = 51
$ws.Range("J15").Value = 46
$ws.Range("J16").Value = 10
$ws.Range("J17").Value = 10
$ws.Range("J18").Value = 60
$ws.Range("J19").Value = 145
$ws.Range("J20").Value = 88
$ws.Range("J29").Value = 237
$ws.Range("J31").Value = 30
$ws.Range("J33").Value = 171
$ws.Range("J37").Value = 146
$ws.Range("J38").Value = 2
$ws.Range("J41").Value = 27
$ws.Range("J42").Value = 178
$ws.Range("J43").Value = 47
$ws.Range("J44").Value = 32
$ws.Range("J46").Value = 13
$ws.Range("J51").Value = 59
$ws.Range("J52").Value = 94
$ws.Range("J53").Value = 39
$ws.Range("J54").Value = 79
$ws.Range("J55").Value = 52
$ws.Range("H63").Value = 237
$ws.Range("J63").Value = 20
$ws.Range("J65").Value = 110
$ws.Range("J67").Value = 149
$ws.Range("J75").Value = 19
$ws.Range("J76").Value = 70
$ws.Range("J77").Value = 34
$ws.Range("J79").Value = 124
$ws.Range("J83").Value = 99
$ws.Range("J84").Value = 42
$ws.Range("J85").Value = 170
$ws.Range("J86").Value = 17
$ws.Range("J87").Value = 17
$ws.Range("J90").Value = 49
$ws.Range("J91").Value = 57
$ws.Range("J92").Value = 13
$ws.Range("I93").Value = 151
$ws.Range("J95").Value = 72
$ws.Range("J96").Value = 57
$ws.Range("H101").Value = 26000
$ws.Range("J101").Value = 4211

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 60
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 23
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 16
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 87
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 43
$ws.Range("I6").Value = 226
$ws.Range("I7").Value = 818
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 48
$ws.Range("J4").Value = 5
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 20
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 38
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 16
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 28
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 13
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 34
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 22
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 20

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 14
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("J2").Value = 1
$ws.Range("J6").Value = 2
